$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "wrote everything" entry on 19 марта to the more specific
#    description of the lockfree-queue benchmark, and add a new day (20 марта)
#    describing the stage-data-type benchmark.
$ws.Range("A25").Value2 = "20 марта"
$ws.Range("B24").Value2 = "Правка теоритической части текста в черновике вкр. Создание проекта для проведения бенчмарков различных библиотек, а также написание бенчмарка для тестирования скорости коннекторов (очередей) (бенчмарк testLockfreeQueues)"
$ws.Range("B25").Value2 = "Написание бенчмарка для тестирования структур данных, в которых хранятся данные, которые обрабатываются стадиями конвейера (бенчмарк testStageDataType)"

# 2) Column B needs to grow to fit the new (longer) text (was bestFit at
#    153.57 chars; the longer benchmark description now needs ~226.57).
$ws.Columns.Item(2).ColumnWidth = 225.667

# 3) Update the active selection / scroll position to match where the user
#    ended up after typing the new row.
$ws.Range("B28").Select() | Out-Null
